$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 77 (shifts existing rows 77-124 down to 78-125)
$ws.Rows.Item(77).Insert()

$ws.Cells.Item(77,1).Value = 10
$ws.Cells.Item(77,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(77,3).Value = "La Araucanía"
$ws.Cells.Item(77,4).Value = 44567
$ws.Cells.Item(77,5).Value = 9
$ws.Cells.Item(77,6).Value = 100114007
$ws.Cells.Item(77,7).Value = "Jengibre"
$ws.Cells.Item(77,8).Value = "Sin especificar"
$ws.Cells.Item(77,9).Value = "Primera"
$ws.Cells.Item(77,10).Value = 50
$ws.Cells.Item(77,11).Value = 20000
$ws.Cells.Item(77,12).Value = 20000
$ws.Cells.Item(77,13).Value = 20000
$ws.Cells.Item(77,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(77,15).Value = "Perú"
$ws.Cells.Item(77,16).Value = 1538
$ws.Cells.Item(77,17).Value = 13
$ws.Cells.Item(77,18).Value = "Hortaliza"

# Insert a second new record row at row 121 (shifts rows 121-125 down to 122-126)
$ws.Rows.Item(121).Insert()

$ws.Cells.Item(121,1).Value = 10
$ws.Cells.Item(121,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(121,3).Value = "La Araucanía"
$ws.Cells.Item(121,4).Value = 44568
$ws.Cells.Item(121,5).Value = 9
$ws.Cells.Item(121,6).Value = 100114007
$ws.Cells.Item(121,7).Value = "Jengibre"
$ws.Cells.Item(121,8).Value = "Sin especificar"
$ws.Cells.Item(121,9).Value = "Primera"
$ws.Cells.Item(121,10).Value = 40
$ws.Cells.Item(121,11).Value = 20000
$ws.Cells.Item(121,12).Value = 20000
$ws.Cells.Item(121,13).Value = 20000
$ws.Cells.Item(121,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(121,15).Value = "Perú"
$ws.Cells.Item(121,16).Value = 1538
$ws.Cells.Item(121,17).Value = 13
$ws.Cells.Item(121,18).Value = "Hortaliza"
